$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append ".age_premiere_conso" to each header cell (B1:P1) in row 1,
# right after the existing ".jamais" suffix.
$cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P")
foreach ($col in $cols) {
    $cell = $ws.Range("$col`1")
    $cell.Value2 = $cell.Value2.ToString() + ".age_premiere_conso"
}
